$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheets involved
# ------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Status column = C on both language sheets)
# ------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ------------------------------------------------------------------
# 2. Fill in "Latest Target File" (I) / "Latest Handback File" (J) /
#    "Latest Handback DateTime" (K) now that the handback completed.
# ------------------------------------------------------------------

# --- zh-cn sheet ---
$zhcn.Range("I2").Value = "365764c5-d128-40bc-9cee-edb6cb33f643.md"
$zhcn.Range("J2").Value = "365764c5-d128-40bc-9cee-edb6cb33f643.1ce20cdbf5cf3afd932102f751cd6f9bea017291.zh-cn.xlf"

$zhcn.Range("I3").Value = "a997d19f-6a67-4018-8d32-d9177a7f1463.md"
$zhcn.Range("J3").Value = "a997d19f-6a67-4018-8d32-d9177a7f1463.68579c22ce767fa74bb61e44336f2d9a51bbe549.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-31 21:17:29"
$zhcn.Range("K3").Value = "2016-08-31 21:17:29"

# --- de-de sheet ---
$dede.Range("I2").Value = "365764c5-d128-40bc-9cee-edb6cb33f643.md"
$dede.Range("J2").Value = "365764c5-d128-40bc-9cee-edb6cb33f643.1ce20cdbf5cf3afd932102f751cd6f9bea017291.de-de.xlf"

$dede.Range("I3").Value = "a997d19f-6a67-4018-8d32-d9177a7f1463.md"
$dede.Range("J3").Value = "a997d19f-6a67-4018-8d32-d9177a7f1463.68579c22ce767fa74bb61e44336f2d9a51bbe549.de-de.xlf"

$dede.Range("K2").Value = "2016-08-31 21:17:38"
$dede.Range("K3").Value = "2016-08-31 21:17:38"

# ------------------------------------------------------------------
# 3. Re-create the hyperlinks on both sheets: keep the existing ones
#    on column A ("...md" source file) and add matching hyperlinks on
#    column I ("Latest Target File"), which now also shows the ".md"
#    file name now that the handback is complete.
# ------------------------------------------------------------------
$target1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/365764c5-d128-40bc-9cee-edb6cb33f643.md"
$target2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/a997d19f-6a67-4018-8d32-d9177a7f1463.md"

foreach ($ws in @($zhcn, $dede)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $target1, "", "", "365764c5-d128-40bc-9cee-edb6cb33f643.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $target1, "", "", "365764c5-d128-40bc-9cee-edb6cb33f643.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $target2, "", "", "a997d19f-6a67-4018-8d32-d9177a7f1463.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $target2, "", "", "a997d19f-6a67-4018-8d32-d9177a7f1463.md")
}

# ------------------------------------------------------------------
# 4. Widen the columns that now show the longer file names.
# ------------------------------------------------------------------

# Overview sheet: zh-cn (E) / de-de (F) status columns
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

foreach ($ws in @($zhcn, $dede)) {
    $ws.Range("C1").ColumnWidth = 29.9777047293527
    $ws.Range("I1").ColumnWidth = 40
    $ws.Range("J1").ColumnWidth = 40
}
